# Update the "as_of_utc" timestamp column (AA) from 2025-11-19 07:08:30
# to 2025-11-19 10:05:30 on both data sheets ("Главные" and "Линейные").

$wb = $excel.ActiveWorkbook

$oldValue = "2025-11-19 07:08:30"
$newValue = "2025-11-19 10:05:30"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Cells.Item($row, 27)  # Column AA = 27
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}
